$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.1190320826869504;  C = 0.306821227259698;    D = 0.7527432677738641; E = 0.4942365360607697; G = 1.672833113781282 }
    3 = @{ B = 0.04271373187048222; C = 0.002571899574220771; D = 0.1494219747398047; E = 10.19245300693656;  G = 10.38716061312106 }
    4 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    5 = @{ B = 0.1190320826869504;  C = 0.04071648406533734;  D = 0.1494219747398047; E = 0.4942365360607697; G = 0.8034070775528621 }
    6 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 3.537761648806719;  E = 0.4942365360607697; G = 8.974608811992548 }
    7 = @{ B = 0.6606524410359556;  C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 3.56341032713086 }
    8 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
